$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> MuSCs
$ws.Cells.Item(2,1).Value2 = "ECs"
$ws.Cells.Item(2,2).Value2 = "Pdgfb"
$ws.Cells.Item(2,3).Value2 = "Art1"
$ws.Cells.Item(2,4).Value2 = "MuSCs"
$ws.Cells.Item(2,5).Value2 = 3
$ws.Cells.Item(2,6).Value2 = 1
$ws.Cells.Item(2,7).Value2 = 47.25342666666666
$ws.Cells.Item(2,8).Value2 = 141.76028
$ws.Cells.Item(2,9).Value2 = 0.7256581002375004
$ws.Cells.Item(2,10).Value2 = 0.7256581002375005
$ws.Cells.Item(2,11).Value2 = 3
$ws.Cells.Item(2,12).Value2 = 1
$ws.Cells.Item(2,13).Value2 = 0.2388073333333333
$ws.Cells.Item(2,14).Value2 = 0.716422
$ws.Cells.Item(2,15).Value2 = 0.9047137371081621
$ws.Cells.Item(2,16).Value2 = 0.9047137371081621
$ws.Cells.Item(2,17).Value2 = 11.28446481312889
$ws.Cells.Item(2,18).Value2 = 101.56018331816
$ws.Cells.Item(2,19).Value2 = 0.6565128517286782
$ws.Cells.Item(2,20).Value2 = 0.6565128517286783

# Row 3: ECs -> Resolving-Mac
$ws.Cells.Item(3,1).Value2 = "ECs"
$ws.Cells.Item(3,2).Value2 = "Pdgfb"
$ws.Cells.Item(3,3).Value2 = "Art1"
$ws.Cells.Item(3,4).Value2 = "Resolving-Mac"
$ws.Cells.Item(3,5).Value2 = 3
$ws.Cells.Item(3,6).Value2 = 1
$ws.Cells.Item(3,7).Value2 = 47.25342666666666
$ws.Cells.Item(3,8).Value2 = 141.76028
$ws.Cells.Item(3,9).Value2 = 0.7256581002375004
$ws.Cells.Item(3,10).Value2 = 0.7256581002375005
$ws.Cells.Item(3,11).Value2 = 1
$ws.Cells.Item(3,12).Value2 = 0.3333333333333333
$ws.Cells.Item(3,13).Value2 = 0.02515166666666667
$ws.Cells.Item(3,14).Value2 = 0.07545499999999999
$ws.Cells.Item(3,15).Value2 = 0.095286262891838
$ws.Cells.Item(3,16).Value2 = 0.095286262891838
$ws.Cells.Item(3,17).Value2 = 1.188502436377778
$ws.Cells.Item(3,18).Value2 = 10.6965219274
$ws.Cells.Item(3,19).Value2 = 0.0691452485088222
$ws.Cells.Item(3,20).Value2 = 0.06914524850882221

# Row 4: FAPs -> MuSCs
$ws.Cells.Item(4,1).Value2 = "FAPs"
$ws.Cells.Item(4,2).Value2 = "Pdgfb"
$ws.Cells.Item(4,3).Value2 = "Art1"
$ws.Cells.Item(4,4).Value2 = "MuSCs"
$ws.Cells.Item(4,5).Value2 = 1
$ws.Cells.Item(4,6).Value2 = 0.3333333333333333
$ws.Cells.Item(4,7).Value2 = 0.05229733333333333
$ws.Cells.Item(4,8).Value2 = 0.156892
$ws.Cells.Item(4,9).Value2 = 0.0008031160114981568
$ws.Cells.Item(4,10).Value2 = 0.0008031160114981569
$ws.Cells.Item(4,11).Value2 = 3
$ws.Cells.Item(4,12).Value2 = 1
$ws.Cells.Item(4,13).Value2 = 0.2388073333333333
$ws.Cells.Item(4,14).Value2 = 0.716422
$ws.Cells.Item(4,15).Value2 = 0.9047137371081621
$ws.Cells.Item(4,16).Value2 = 0.9047137371081621
$ws.Cells.Item(4,17).Value2 = 0.01248898671377778
$ws.Cells.Item(4,18).Value2 = 0.112400880424
$ws.Cells.Item(4,19).Value2 = 0.0007265900880938991
$ws.Cells.Item(4,20).Value2 = 0.0007265900880938992

# Row 5: FAPs -> Resolving-Mac
$ws.Cells.Item(5,1).Value2 = "FAPs"
$ws.Cells.Item(5,2).Value2 = "Pdgfb"
$ws.Cells.Item(5,3).Value2 = "Art1"
$ws.Cells.Item(5,4).Value2 = "Resolving-Mac"
$ws.Cells.Item(5,5).Value2 = 1
$ws.Cells.Item(5,6).Value2 = 0.3333333333333333
$ws.Cells.Item(5,7).Value2 = 0.05229733333333333
$ws.Cells.Item(5,8).Value2 = 0.156892
$ws.Cells.Item(5,9).Value2 = 0.0008031160114981568
$ws.Cells.Item(5,10).Value2 = 0.0008031160114981569
$ws.Cells.Item(5,11).Value2 = 1
$ws.Cells.Item(5,12).Value2 = 0.3333333333333333
$ws.Cells.Item(5,13).Value2 = 0.02515166666666667
$ws.Cells.Item(5,14).Value2 = 0.07545499999999999
$ws.Cells.Item(5,15).Value2 = 0.095286262891838
$ws.Cells.Item(5,16).Value2 = 0.095286262891838
$ws.Cells.Item(5,17).Value2 = 0.001315365095555556
$ws.Cells.Item(5,18).Value2 = 0.01183828586
$ws.Cells.Item(5,19).Value2 = 0.00007652592340425776
$ws.Cells.Item(5,20).Value2 = 0.00007652592340425776

# Row 6: Inflammatory-Mac -> MuSCs
$ws.Cells.Item(6,1).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(6,2).Value2 = "Pdgfb"
$ws.Cells.Item(6,3).Value2 = "Art1"
$ws.Cells.Item(6,4).Value2 = "MuSCs"
$ws.Cells.Item(6,5).Value2 = 3
$ws.Cells.Item(6,6).Value2 = 1
$ws.Cells.Item(6,7).Value2 = 6.996562666666667
$ws.Cells.Item(6,8).Value2 = 20.989688
$ws.Cells.Item(6,9).Value2 = 0.1074443216298519
$ws.Cells.Item(6,10).Value2 = 0.1074443216298519
$ws.Cells.Item(6,11).Value2 = 3
$ws.Cells.Item(6,12).Value2 = 1
$ws.Cells.Item(6,13).Value2 = 0.2388073333333333
$ws.Cells.Item(6,14).Value2 = 0.716422
$ws.Cells.Item(6,15).Value2 = 0.9047137371081621
$ws.Cells.Item(6,16).Value2 = 0.9047137371081621
$ws.Cells.Item(6,17).Value2 = 1.670830472926222
$ws.Cells.Item(6,18).Value2 = 15.037474256336
$ws.Cells.Item(6,19).Value2 = 0.09720635375279464
$ws.Cells.Item(6,20).Value2 = 0.09720635375279466

# Row 7: Inflammatory-Mac -> Resolving-Mac
$ws.Cells.Item(7,1).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(7,2).Value2 = "Pdgfb"
$ws.Cells.Item(7,3).Value2 = "Art1"
$ws.Cells.Item(7,4).Value2 = "Resolving-Mac"
$ws.Cells.Item(7,5).Value2 = 3
$ws.Cells.Item(7,6).Value2 = 1
$ws.Cells.Item(7,7).Value2 = 6.996562666666667
$ws.Cells.Item(7,8).Value2 = 20.989688
$ws.Cells.Item(7,9).Value2 = 0.1074443216298519
$ws.Cells.Item(7,10).Value2 = 0.1074443216298519
$ws.Cells.Item(7,11).Value2 = 1
$ws.Cells.Item(7,12).Value2 = 0.3333333333333333
$ws.Cells.Item(7,13).Value2 = 0.02515166666666667
$ws.Cells.Item(7,14).Value2 = 0.07545499999999999
$ws.Cells.Item(7,15).Value2 = 0.095286262891838
$ws.Cells.Item(7,16).Value2 = 0.095286262891838
$ws.Cells.Item(7,17).Value2 = 0.1759752120044444
$ws.Cells.Item(7,18).Value2 = 1.58377690804
$ws.Cells.Item(7,19).Value2 = 0.01023796787705726
$ws.Cells.Item(7,20).Value2 = 0.01023796787705727

# Row 8: MuSCs -> MuSCs
$ws.Cells.Item(8,1).Value2 = "MuSCs"
$ws.Cells.Item(8,2).Value2 = "Pdgfb"
$ws.Cells.Item(8,3).Value2 = "Art1"
$ws.Cells.Item(8,4).Value2 = "MuSCs"
$ws.Cells.Item(8,5).Value2 = 3
$ws.Cells.Item(8,6).Value2 = 1
$ws.Cells.Item(8,7).Value2 = 2.605199666666667
$ws.Cells.Item(8,8).Value2 = 7.815599000000001
$ws.Cells.Item(8,9).Value2 = 0.0400073470689964
$ws.Cells.Item(8,10).Value2 = 0.0400073470689964
$ws.Cells.Item(8,11).Value2 = 3
$ws.Cells.Item(8,12).Value2 = 1
$ws.Cells.Item(8,13).Value2 = 0.2388073333333333
$ws.Cells.Item(8,14).Value2 = 0.716422
$ws.Cells.Item(8,15).Value2 = 0.9047137371081621
$ws.Cells.Item(8,16).Value2 = 0.9047137371081621
$ws.Cells.Item(8,17).Value2 = 0.6221407851975557
$ws.Cells.Item(8,18).Value2 = 5.599267066778
$ws.Cells.Item(8,19).Value2 = 0.03619519647857501
$ws.Cells.Item(8,20).Value2 = 0.03619519647857501

# Row 9: MuSCs -> Resolving-Mac
$ws.Cells.Item(9,1).Value2 = "MuSCs"
$ws.Cells.Item(9,2).Value2 = "Pdgfb"
$ws.Cells.Item(9,3).Value2 = "Art1"
$ws.Cells.Item(9,4).Value2 = "Resolving-Mac"
$ws.Cells.Item(9,5).Value2 = 3
$ws.Cells.Item(9,6).Value2 = 1
$ws.Cells.Item(9,7).Value2 = 2.605199666666667
$ws.Cells.Item(9,8).Value2 = 7.815599000000001
$ws.Cells.Item(9,9).Value2 = 0.0400073470689964
$ws.Cells.Item(9,10).Value2 = 0.0400073470689964
$ws.Cells.Item(9,11).Value2 = 1
$ws.Cells.Item(9,12).Value2 = 0.3333333333333333
$ws.Cells.Item(9,13).Value2 = 0.02515166666666667
$ws.Cells.Item(9,14).Value2 = 0.07545499999999999
$ws.Cells.Item(9,15).Value2 = 0.095286262891838
$ws.Cells.Item(9,16).Value2 = 0.095286262891838
$ws.Cells.Item(9,17).Value2 = 0.06552511361611112
$ws.Cells.Item(9,18).Value2 = 0.589726022545
$ws.Cells.Item(9,19).Value2 = 0.003812150590421396
$ws.Cells.Item(9,20).Value2 = 0.003812150590421396

# Row 10: Resolving-Mac -> MuSCs
$ws.Cells.Item(10,1).Value2 = "Resolving-Mac"
$ws.Cells.Item(10,2).Value2 = "Pdgfb"
$ws.Cells.Item(10,3).Value2 = "Art1"
$ws.Cells.Item(10,4).Value2 = "MuSCs"
$ws.Cells.Item(10,5).Value2 = 3
$ws.Cells.Item(10,6).Value2 = 1
$ws.Cells.Item(10,7).Value2 = 8.210544666666665
$ws.Cells.Item(10,8).Value2 = 24.631634
$ws.Cells.Item(10,9).Value2 = 0.126087115052153
$ws.Cells.Item(10,10).Value2 = 0.126087115052153
$ws.Cells.Item(10,11).Value2 = 3
$ws.Cells.Item(10,12).Value2 = 1
$ws.Cells.Item(10,13).Value2 = 0.2388073333333333
$ws.Cells.Item(10,14).Value2 = 0.716422
$ws.Cells.Item(10,15).Value2 = 0.9047137371081621
$ws.Cells.Item(10,16).Value2 = 0.9047137371081621
$ws.Cells.Item(10,17).Value2 = 1.960738277060889
$ws.Cells.Item(10,18).Value2 = 17.646644493548
$ws.Cells.Item(10,19).Value2 = 0.1140727450600201
$ws.Cells.Item(10,20).Value2 = 0.1140727450600202

# Row 11: Resolving-Mac -> Resolving-Mac
$ws.Cells.Item(11,1).Value2 = "Resolving-Mac"
$ws.Cells.Item(11,2).Value2 = "Pdgfb"
$ws.Cells.Item(11,3).Value2 = "Art1"
$ws.Cells.Item(11,4).Value2 = "Resolving-Mac"
$ws.Cells.Item(11,5).Value2 = 3
$ws.Cells.Item(11,6).Value2 = 1
$ws.Cells.Item(11,7).Value2 = 8.210544666666665
$ws.Cells.Item(11,8).Value2 = 24.631634
$ws.Cells.Item(11,9).Value2 = 0.126087115052153
$ws.Cells.Item(11,10).Value2 = 0.126087115052153
$ws.Cells.Item(11,11).Value2 = 1
$ws.Cells.Item(11,12).Value2 = 0.3333333333333333
$ws.Cells.Item(11,13).Value2 = 0.02515166666666667
$ws.Cells.Item(11,14).Value2 = 0.07545499999999999
$ws.Cells.Item(11,15).Value2 = 0.095286262891838
$ws.Cells.Item(11,16).Value2 = 0.095286262891838
$ws.Cells.Item(11,17).Value2 = 0.2065088826077777
$ws.Cells.Item(11,18).Value2 = 1.85857994347
$ws.Cells.Item(11,19).Value2 = 0.01201436999213288
$ws.Cells.Item(11,20).Value2 = 0.01201436999213288
